$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" column cells that look like plain numbers (e.g. "215.96") must stay
# TEXT cells (as in the source file, t="inlineStr"), matching values such as
# "26.024.42" that cannot be numeric. Assigning such a string straight to
# .Value lets Excel auto-coerce it to a float (losing the "." thousands-style
# formatting / introducing FP rounding noise), so instead we stage the text in
# a scratch cell formatted as Text and PasteSpecial only the value across - this
# keeps the destination cell a string without altering its number format/style.
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
$helper = $ws.Range("Z1")

$ws.Range('D2').Value = '26.010.12'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '1.642.62'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.41%  '
$helper.NumberFormat = '@'
$helper.Value = '215.96'
$helper.Copy()
$ws.Range('D5').PasteSpecial($xlPasteValues)
$ws.Range('E5').Value = '  +0.78%  '
$helper.NumberFormat = '@'
$helper.Value = '0.507'
$helper.Copy()
$ws.Range('D6').PasteSpecial($xlPasteValues)
$ws.Range('E6').Value = '  +0.66%  '
$helper.NumberFormat = '@'
$helper.Value = '1.01'
$helper.Copy()
$ws.Range('D7').PasteSpecial($xlPasteValues)
$ws.Range('E7').Value = '  +0.39%  '
$helper.NumberFormat = '@'
$helper.Value = '0.255'
$helper.Copy()
$ws.Range('D8').PasteSpecial($xlPasteValues)
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  +1.19%  '
$helper.NumberFormat = '@'
$helper.Value = '19.61'
$helper.Copy()
$ws.Range('D10').PasteSpecial($xlPasteValues)
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '1.682.42'
$ws.Range('E12').Value = '  +2.69%  '
$ws.Range('D13').Value = '1.870.82'
$ws.Range('E13').Value = '  +0.89%  '
$helper.NumberFormat = '@'
$helper.Value = '4.27'
$helper.Copy()
$ws.Range('D14').PasteSpecial($xlPasteValues)
$ws.Range('E14').Value = '  +0.64%  '
$helper.NumberFormat = '@'
$helper.Value = '0.543'
$helper.Copy()
$ws.Range('D15').PasteSpecial($xlPasteValues)
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('E16').Value = '  +1.02%  '
$helper.NumberFormat = '@'
$helper.Value = '63.44'
$helper.Copy()
$ws.Range('D17').PasteSpecial($xlPasteValues)
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').Value = '26.093.85'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  +0.45%  '
$helper.NumberFormat = '@'
$helper.Value = '195.15'
$helper.Copy()
$ws.Range('D20').PasteSpecial($xlPasteValues)
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('E21').Value = '  -0.45%  '
$helper.NumberFormat = '@'
$helper.Value = '9.92'
$helper.Copy()
$ws.Range('D22').PasteSpecial($xlPasteValues)
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$helper.NumberFormat = '@'
$helper.Value = '1.79'
$helper.Copy()
$ws.Range('D24').PasteSpecial($xlPasteValues)
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$helper.NumberFormat = '@'
$helper.Value = '0.131'
$helper.Copy()
$ws.Range('D25').PasteSpecial($xlPasteValues)
$ws.Range('E25').Value = '  +4.08%  '
$ws.Range('E26').Value = '  +0.45%  '
$helper.NumberFormat = '@'
$helper.Value = '143.17'
$helper.Copy()
$ws.Range('D27').PasteSpecial($xlPasteValues)
$ws.Range('E27').Value = '  +0.14%  '
$helper.NumberFormat = '@'
$helper.Value = '6.88'
$helper.Copy()
$ws.Range('D28').PasteSpecial($xlPasteValues)
$ws.Range('E28').Value = '  +0.93%  '
$helper.NumberFormat = '@'
$helper.Value = '15.55'
$helper.Copy()
$ws.Range('D29').PasteSpecial($xlPasteValues)
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('E35').Value = '  +1.89%  '
$helper.NumberFormat = '@'
$helper.Value = '0.905'
$helper.Copy()
$ws.Range('D36').PasteSpecial($xlPasteValues)
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').Value = '1.129.76'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('E38').Value = '  -1.32%  '
$ws.Range('E39').Value = '  -0.07%  '
$helper.NumberFormat = '@'
$helper.Value = '0.0156'
$helper.Copy()
$ws.Range('D40').PasteSpecial($xlPasteValues)
$ws.Range('E40').Value = '  +0.66%  '
$helper.NumberFormat = '@'
$helper.Value = '5.47'
$helper.Copy()
$ws.Range('D41').PasteSpecial($xlPasteValues)
$ws.Range('E41').Value = '  +0.76%  '
$helper.NumberFormat = '@'
$helper.Value = '99.30'
$helper.Copy()
$ws.Range('D42').PasteSpecial($xlPasteValues)
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('D44').Value = '1.780.46'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('E45').Value = '  +4.54%  '
$helper.NumberFormat = '@'
$helper.Value = '56.71'
$helper.Copy()
$ws.Range('D46').PasteSpecial($xlPasteValues)
$ws.Range('E46').Value = '  +1.23%  '
$helper.NumberFormat = '@'
$helper.Value = '0.0523'
$helper.Copy()
$ws.Range('D47').PasteSpecial($xlPasteValues)
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('E48').Value = '  +0.97%  '
$helper.NumberFormat = '@'
$helper.Value = '7.77'
$helper.Copy()
$ws.Range('D49').PasteSpecial($xlPasteValues)
$ws.Range('E49').Value = '  +2.35%  '
$ws.Range('E50').Value = '  -0.12%  '
$helper.NumberFormat = '@'
$helper.Value = '0.0954'
$helper.Copy()
$ws.Range('D51').PasteSpecial($xlPasteValues)
$ws.Range('E51').Value = '  -0.39%  '

$helper.Clear()
$excel.CutCopyMode = $false
